$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New label near the first ("old way") lookup table: L2 = "Old way"
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "Old way"

# ---------------------------------------------------------------------------
# 1b. "New way" heading + r[row][col] caption (row 38) -- the shared strings
#      for these are interned right after "Old way" and before the rest of
#      the lookup-table contents.
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "New way"
$ws.Range("D38").Value = "r[row][col]"

# ---------------------------------------------------------------------------
# 2. "Old way" r[row][col] lookup table (rows 31-37)
# ---------------------------------------------------------------------------
$ws.Range("J31").Value = "r0"
$ws.Range("K31").Value = "r1"
$ws.Range("L31").Value = "r2"
$ws.Range("M31").Value = "r3"
$ws.Range("N31").Value = "r4"
$ws.Range("O31").Value = "r5"
$ws.Range("P31").Value = "r6"
$ws.Range("Q31").Value = "r7"

$ws.Range("J32").Value = "c0"
$ws.Range("K32").Value = "c1"
$ws.Range("L32").Value = "c2"
$ws.Range("M32").Value = "c3"
$ws.Range("N32").Value = "c4"
$ws.Range("O32").Value = "c5"
$ws.Range("P32").Value = "c6"
$ws.Range("Q32").Value = "c7"

$ws.Range("J33").Value = "c8"
$ws.Range("J34").Value = "c16"
$ws.Range("J35").Value = "c24"
$ws.Range("O34").Value = "c21"
$ws.Range("J36").Value = "c32"
$ws.Range("J37").Value = "c40"

$ws.Range("K33").Value = "c9"
$ws.Range("L33").Value = "c10"
$ws.Range("M33").Value = "c11"
$ws.Range("N33").Value = "c12"
$ws.Range("O33").Value = "c13"
$ws.Range("P33").Value = "c14"
$ws.Range("Q33").Value = "c15"

$ws.Range("K34").Value = "c17"
$ws.Range("L34").Value = "c18"
$ws.Range("M34").Value = "c19"
$ws.Range("N34").Value = "c20"
$ws.Range("P34").Value = "c22"
$ws.Range("Q34").Value = "c23"

$ws.Range("K35").Value = "c25"
$ws.Range("L35").Value = "c26"
$ws.Range("M35").Value = "c27"
$ws.Range("N35").Value = "c28"
$ws.Range("O35").Value = "c29"
$ws.Range("P35").Value = "c30"
$ws.Range("Q35").Value = "c31"

$ws.Range("K36").Value = "c33"
$ws.Range("L36").Value = "c34"
$ws.Range("M36").Value = "c35"
$ws.Range("N36").Value = "c36"
$ws.Range("O36").Value = "c37"
$ws.Range("P36").Value = "c38"
$ws.Range("Q36").Value = "c39"

$ws.Range("K37").Value = "c41"
$ws.Range("L37").Value = "c42"
$ws.Range("M37").Value = "c43"
$ws.Range("N37").Value = "c44"
$ws.Range("O37").Value = "c45"
$ws.Range("P37").Value = "c46"
$ws.Range("Q37").Value = "c47"

$ws.Range("S32").Value = "1 left"
$ws.Range("S33").Value = "2 left"
$ws.Range("S34").Value = "3 left"
$ws.Range("S35").Value = "4 left"
$ws.Range("S31").Value = "0 left"
$ws.Range("S36").Value = "5 left"
$ws.Range("S37").Value = "7 left"

# ---------------------------------------------------------------------------
# 4. "New way" color-coded, merged lookup grid (rows 41,43,45,47,51,53,55,57)
# ---------------------------------------------------------------------------
$blockCols = @(
    @{Start="B";  End="I";  Fmt="BD3"},
    @{Start="J";  End="Q";  Fmt="AF26"},
    @{Start="R";  End="Y";  Fmt="AF3"},
    @{Start="Z";  End="AG"; Fmt="AV3"},
    @{Start="AH"; End="AO"; Fmt="BL3"},
    @{Start="AP"; End="AW"; Fmt="AN3"},
    @{Start="AX"; End="BE"; Fmt="BT3"},
    @{Start="BF"; End="BM"; Fmt=$null}
)

$gridRows = @(41, 43, 45, 47, 51, 53, 55, 57)
$gridValues = @{
    41 = @("r0","r1","r2","r3","r4","r5","r6","r7")
    43 = @("c0","c1","c2","c3","c4","c5","c6","c7")
    45 = @("c8","c9","c10","c11","c12","c13","c14","c15")
    47 = @("c16","c17","c18","c19","c20","c21","c22","c23")
    51 = @("c24","c25","c26","c27","c28","c29","c30","c31")
    53 = @("c32","c33","c34","c35","c36","c37","c38","c39")
    55 = @("c40","c41","c42","c43","c44","c45","c46","c47")
    57 = @($null,$null,$null,$null,$null,$null,$null,$null)
}

foreach ($r in $gridRows) {
    $vals = $gridValues[$r]
    for ($i = 0; $i -lt $blockCols.Length; $i++) {
        $blk = $blockCols[$i]
        $rng = $ws.Range("$($blk.Start)$r`:$($blk.End)$r")
        if ($blk.Fmt) {
            $ws.Range($blk.Fmt).Copy()
            $rng.PasteSpecial(-4122)
        } else {
            $rng.Interior.Color = 65535
        }
        $rng.HorizontalAlignment = -4108
        $rng.Merge()
        $v = $vals[$i]
        if ($v) {
            $rng.Value = $v
        }
    }
}

# ---------------------------------------------------------------------------
# 5. Sheet view: selection on J41:Q41, scrolled so row 24 is at the top
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J41:Q41").Select()
